$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update Maximo (C2) ---
$ws1 = $wb.Worksheets.Item("Resumen")
$ws1.Range("C2").Value = 590.9590496269864

# --- Sheet "Solucion": reorder Salida assignments (column B) ---
$ws2 = $wb.Worksheets.Item("Solucion")
$ws2.Range("B2").Value = "S033"
$ws2.Range("B5").Value = "S039"
$ws2.Range("B7").Value = "S001"
$ws2.Range("B8").Value = "S034"
$ws2.Range("B9").Value = "S044"
$ws2.Range("B10").Value = "S007"
$ws2.Range("B11").Value = "S036"
$ws2.Range("B12").Value = "S030"
$ws2.Range("B13").Value = "S043"
$ws2.Range("B14").Value = "S046"
$ws2.Range("B15").Value = "S031"
$ws2.Range("B17").Value = "S026"
$ws2.Range("B18").Value = "S037"
$ws2.Range("B19").Value = "S045"
$ws2.Range("B20").Value = "S005"
$ws2.Range("B21").Value = "S040"
$ws2.Range("B22").Value = "S003"
$ws2.Range("B23").Value = "S028"
$ws2.Range("B24").Value = "S016"
$ws2.Range("B25").Value = "S032"
$ws2.Range("B26").Value = "S006"
$ws2.Range("B28").Value = "S050"
$ws2.Range("B29").Value = "S048"
$ws2.Range("B30").Value = "S004"
$ws2.Range("B31").Value = "S035"
$ws2.Range("B32").Value = "S047"
$ws2.Range("B33").Value = "S049"
$ws2.Range("B34").Value = "S013"
$ws2.Range("B35").Value = "S054"
$ws2.Range("B36").Value = "S014"
$ws2.Range("B37").Value = "S008"
$ws2.Range("B38").Value = "S038"
$ws2.Range("B39").Value = "S042"
$ws2.Range("B40").Value = "S010"
$ws2.Range("B41").Value = "S021"
$ws2.Range("B42").Value = "S055"
$ws2.Range("B43").Value = "S011"
$ws2.Range("B44").Value = "S052"
$ws2.Range("B45").Value = "S015"
$ws2.Range("B46").Value = "S056"
$ws2.Range("B47").Value = "S051"
$ws2.Range("B48").Value = "S009"
$ws2.Range("B49").Value = "S053"
$ws2.Range("B51").Value = "S017"
$ws2.Range("B52").Value = "S062"
$ws2.Range("B53").Value = "S061"
$ws2.Range("B55").Value = "S012"
$ws2.Range("B56").Value = "S018"
$ws2.Range("B59").Value = "S063"
$ws2.Range("B60").Value = "S019"

# --- Sheet "Metricas": update Tiempo values (column B) ---
$ws3 = $wb.Worksheets.Item("Metricas")
$ws3.Range("B2").Value = 590.9590496269864
$ws3.Range("B3").Value = 579.0398151151477
$ws3.Range("B4").Value = 590.6868580387069
